$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 23:22"

# Estados Unidos (row 4): refreshed totals
$ws.Range("B4").Value = 1029878
$ws.Range("C4").Value = 19522
$ws.Range("E4").Value = 831100
$ws.Range("G4").Value = 1843
$ws.Range("H4").Value = 58640

# Francia (row 7): refreshed totals
$ws.Range("B7").Value = 165911
$ws.Range("C7").Value = 2638
$ws.Range("E7").Value = 96738

# Peru overtakes Suiza in the ranking -> rows 19/20 swap identities
# Row 19 becomes Peru with refreshed totals
$ws.Range("A19").Value = "Peru"
$ws.Range("B19").Value = 31190
$ws.Range("C19").Value = 2491
$ws.Range("D19").Value = 8425
$ws.Range("E19").Value = 21911
$ws.Range("F19").Value = 598
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 854

# Row 20 becomes Suiza, carrying the previous Suiza totals
$ws.Range("A20").Value = "Suiza"
$ws.Range("B20").Value = 29264
$ws.Range("C20").Value = 100
$ws.Range("D20").Value = 22600
$ws.Range("E20").Value = 4965
$ws.Range("F20").Value = 185
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 1699
